$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the TPM-derived values in row 2 (M2:T2)
$ws.Range("M2").Value = 0.174999
$ws.Range("N2").Value = 0.349998
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.21229887019
$ws.Range("R2").Value = 1.27379322114
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Remove row 3 (the Neutrophils target-cluster row) entirely
$ws.Rows.Item(3).Delete()
